$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.032.01"
$ws.Range("E2").Value = "  -7.32%  "
$ws.Range("D3").Value = "3.519.74"
$ws.Range("E3").Value = "  -2.82%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "388.95"
$ws.Range("E5").Value = "  -7.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "121.36"
$ws.Range("E6").Value = "  -6.45%  "
$ws.Range("D7").Value = "3.508.57"
$ws.Range("E7").Value = "  -2.87%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.579"
$ws.Range("E8").Value = "  -11.28%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.670"
$ws.Range("E10").Value = "  -12.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.149"
$ws.Range("E11").Value = "  -23.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000312"
$ws.Range("E12").Value = "  -28.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.32"
$ws.Range("E13").Value = "  -8.51%  "
$ws.Range("D14").Value = "4.074.82"
$ws.Range("E14").Value = "  -2.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.04"
$ws.Range("E15").Value = "  -7.77%  "
$ws.Range("E16").Value = "  -3.15%  "
$ws.Range("D17").Value = "3.514.43"
$ws.Range("E17").Value = "  -2.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.96"
$ws.Range("E18").Value = "  +5.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.50"
$ws.Range("E19").Value = "  -7.59%  "
$ws.Range("D20").Value = "63.047.17"
$ws.Range("E20").Value = "  -7.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  -9.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "388.28"
$ws.Range("E22").Value = "  -15.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.81"
$ws.Range("E23").Value = "  +2.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.89"
$ws.Range("E24").Value = "  -8.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.85"
$ws.Range("E25").Value = "  -6.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.45"
$ws.Range("E26").Value = "  +10.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.24"
$ws.Range("E27").Value = "  -6.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.94"
$ws.Range("E28").Value = "  -9.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.54"
$ws.Range("E29").Value = "  -15.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.63"
$ws.Range("E30").Value = "  -5.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.61"
$ws.Range("E31").Value = "  -4.90%  "
$ws.Range("E32").Value = "  -7.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.63"
$ws.Range("E33").Value = "  -8.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.146"
$ws.Range("E35").Value = "  -6.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "36.15"
$ws.Range("E36").Value = "  -9.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.84"
$ws.Range("E37").Value = "  -5.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0430"
$ws.Range("E38").Value = "  -12.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.996"
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.65"
$ws.Range("E40").Value = "  +1.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.129"
$ws.Range("E41").Value = "  -13.55%  "
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").Value = "0.0₃0616"
$ws.Range("E42").Value = "  -22.70%  "
$ws.Range("E43").Value = "  +14.65%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.43"
$ws.Range("E44").Value = "  +20.26%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "138.36"
$ws.Range("E45").Value = "  -6.52%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.93"
$ws.Range("E46").Value = "  -0.24%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.04"
$ws.Range("E47").Value = "  -4.97%  "
$ws.Range("B48").Value = "LidoDAOToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.02"
$ws.Range("E48").Value = "  -6.37%  "
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.44"
$ws.Range("E49").Value = "  -10.46%  "
$ws.Range("E50").Value = "  -10.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.271"
$ws.Range("E51").Value = "  -9.94%  "
